# "update summary bill for stage 01"
# Applies the content edits reflected in the XML diff for bills/summary.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 (Phích cắm điện): quantity 2 -> 3, unit price blank -> 12000
#     (total recalculates automatically: 3*12000 = 36000)
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 12000

# --- Row 8 (order note expanded, price raised 120000 -> 150000)
$ws.Range("B8").Value = "Đơn mua lẻ ở An Lạc + test oscil"
$ws.Range("D8").Value = 150000

# --- Row 7 note: "An" -> "An+Tâm"
$ws.Range("G7").Value = "An+Tâm"

# --- Row 10 label: "Tổng" -> "Thực tổng" (E10 = SUM(E2:E9) recalculates automatically)
$ws.Range("A10").Value = "Thực tổng"

# --- Row 11 (previously a blank spacer row) now carries the actual amount
#     paid ("Tổng "). Give A11:D11 the same bold/merged "total row" look as
#     row 10: merge first, then copy row 10's formatting onto it.
$ws.Range("A11:D11").Merge()
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A11").Value = "Tổng "
$ws.Rows("11:11").RowHeight = 18

$ws.Range("E11").Value = 1800000

# --- Row 12: show the gap between what was actually paid and the real total.
$ws.Range("E12").Formula = "=E11-E10"

$wb.Save()
